$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 22 (week day counter "33"): the Friday/H column cell gets marked
# the same way the rows above it already are (matching the green
# fill + thin border used by H8:H21), and the day-count in column I
# goes from 4 to 5.
$ws.Range("H21").Copy()
$ws.Range("H22").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("I22").Value = 5

# K3 (=125-L3), L3 (=SUM(I3:I28)) and M3 (=L3/125) are formulas and
# will recalculate automatically from the I22 change above.
